# Add two new columns to the worksheet: I ("I0") and J ("IF").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the new header
# cells so they pick up the same bold/border/alignment style used by the
# rest of row 1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J83 (column I = "I0", column J = "IF"), row-by-row.
$iVals = @(8,6,7,5,7,9,6,7,7,8,8,8,7,8,9,7,8,7,8,7,9,7,5,8,9,7,7,7,7,8,8,7,7,6,6,6,8,7,7,8,7,7,7,7,9,7,7,8,9,8,8,6,8,2,5,9,7,8,7,8,6,7,7,7,7,4,6,7,6,8,4,8,6,5,6,4,8,7,6,4,7,6)
$jVals = @(8,6,7,5,8,9,6,7,7,8,8,8,7,8,9,8,8,7,8,7,9,7,5,8,9,7,7,7,7,8,9,7,7,6,6,6,8,7,7,8,7,7,7,7,9,7,7,9,9,9,8,7,8,3,5,9,7,8,7,9,6,8,7,8,7,4,6,7,6,8,4,8,6,6,6,4,8,7,6,4,7,6)

for ($r = 2; $r -le 83; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
